# Allows the drag and drop of CSV files to the attribute table view.
# - Delete the now-unused "Attributes of an element" sheet.
# - Rename the "microm" unit suffix to the proper µm (micro sign) on the
#   remaining "Attributes of a measure" sheet.
# - Move the selection to B7 (from A11).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Attributes of a measure")

# Replace "microm" with "µm" across the surviving sheet's used range.
$micro = [char]0x00B5
$used = $ws1.UsedRange
$used.Replace("microm", $micro + "m") | Out-Null

# Remove the second sheet entirely.
$ws2 = $wb.Worksheets.Item("Attributes of an element")
$ws2.Delete()

# Move the active selection on the remaining sheet to B7.
$ws1.Select()
$ws1.Range("B7").Select()
